# Mark levels 7-1, 7-2, 8-1, 8-2, 9-1, 9-2 (and their bonus rows) as done.
# Columns H, I, J correspond to World 7, World 8, World 9.
# Rows 3 (Graphics) and 4 (Music) are bonus categories; rows 5 (level 1) and
# 6 (level 2) are the first two numbered levels. All of them flip from
# "todo"/"test" to "ok" for these three worlds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3:J6").Value = "ok"

# Leave the selection where the last edit happened, matching the saved file.
$ws.Range("J7").Select()
